# Insert a new data row at row 339 (pushes existing rows 339-448 down to 340-449)
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(339).Insert()

$ws.Range("A339").Value = 5
$ws.Range("B339").Value = "Macroferia Regional de Talca"
$ws.Range("C339").Value = "Maule"
$ws.Range("D339").Value = 44876
$ws.Range("E339").Value = 7
$ws.Range("F339").Value = 100112032
$ws.Range("G339").Value = "Zapallo italiano"
$ws.Range("H339").Value = "Sin especificar"
$ws.Range("I339").Value = "Primera"
$ws.Range("J339").Value = 400
$ws.Range("K339").Value = 10000
$ws.Range("L339").Value = 10000
$ws.Range("M339").Value = 10000
$ws.Range("N339").Value = "`$/caja 50 unidades"
$ws.Range("O339").Value = "Región del Maule"
$ws.Range("P339").Value = 200
$ws.Range("Q339").Value = 50
$ws.Range("R339").Value = "Hortaliza"
